# "Plotnikov Vladimir" is no longer an executor on this task - replace him
# with "Zenkin Nikita" (commit: "Add files via upload").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Зенькин Никита"

# The cursor/selection in the saved file moved from H5 to F7.
$ws.Range("F7").Select()
